# Update the "panel_query_time" style timestamps on the "data" sheet (F2:F11)
# to reflect the later re-run time recorded in the commit, then add a new
# "metadata" worksheet summarizing the panel query.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$times = @(
    "2021-10-05 14:35:26.566529",
    "2021-10-05 14:35:26.566537",
    "2021-10-05 14:35:26.566540",
    "2021-10-05 14:35:26.566543",
    "2021-10-05 14:35:26.566546",
    "2021-10-05 14:35:26.566548",
    "2021-10-05 14:35:26.566551",
    "2021-10-05 14:35:26.566553",
    "2021-10-05 14:35:26.566556",
    "2021-10-05 14:35:26.566558"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}

# Add the new "metadata" worksheet positioned right after "data".
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 2).Value = "Pseudohypoparathyroidism and Albright Hereditary Osteodystrophy"
$meta.Cells.Item(2, 3).Value = 161
# data_version is the literal string "0.4" (not the number 0.4), so force a
# text format before assigning the value to stop Excel auto-converting it.
$meta.Cells.Item(2, 4).NumberFormat = "@"
$meta.Cells.Item(2, 4).Value = "0.4"
$meta.Cells.Item(2, 5).Value = "2020-02-06T06:16:22.371307Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:35:26.562842"
$meta.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/161/?format=json"

# Reuse the header/index style already present on the "data" sheet (bold,
# bordered, centered) for the metadata header row and the leading index cell.
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

# Keep "data" as the active sheet/selection, as it was before this edit.
$ws.Activate()
$ws.Select()
